$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2
$ws.Range("C2").Value = 0.5571428571428572
$ws.Range("J2").Value = 0.01785714285714286
$ws.Range("P2").Value = 0.1642857142857143
$ws.Range("S2").Value = 0.06071428571428571
$ws.Range("B3").Value = 0.02339181286549707
$ws.Range("C3").Value = 0.03508771929824561
$ws.Range("J3").Value = 0.05263157894736842
$ws.Range("P3").Value = 0.7368421052631579
$ws.Range("S3").Value = 0.152046783625731
$ws.Range("J4").Value = 0.08108108108108109
$ws.Range("P4").Value = 0.6216216216216216
$ws.Range("S4").Value = 0.2972972972972973
$ws.Range("B6").Value = 0.07804878048780488
$ws.Range("D6").Value = 0.00975609756097561
$ws.Range("F6").Value = 0.05853658536585366
$ws.Range("J6").Value = 0.2634146341463415
$ws.Range("O6").Value = 0.01463414634146342
$ws.Range("Q6").Value = 0.1219512195121951
$ws.Range("R6").Value = 0.08292682926829269
$ws.Range("S6").Value = 0.3707317073170732
$ws.Range("B7").Value = 0.1379310344827586
$ws.Range("D7").Value = 0.01724137931034483
$ws.Range("F7").Value = 0.05172413793103448
$ws.Range("J7").Value = 0.1379310344827586
$ws.Range("O7").Value = 0.005747126436781609
$ws.Range("Q7").Value = 0.1436781609195402
$ws.Range("R7").Value = 0.07471264367816093
$ws.Range("S7").Value = 0.4310344827586207
$ws.Range("B8").Value = 0.1180722891566265
$ws.Range("D8").Value = 0.01204819277108434
$ws.Range("E8").Value = 0.002409638554216868
$ws.Range("F8").Value = 0.04337349397590361
$ws.Range("J8").Value = 0.108433734939759
$ws.Range("O8").Value = 0.01204819277108434
$ws.Range("Q8").Value = 0.2072289156626506
$ws.Range("R8").Value = 0.0819277108433735
$ws.Range("S8").Value = 0.4144578313253012
$ws.Range("B9").Value = 0.1260504201680672
$ws.Range("D9").Value = 0.01680672268907563
$ws.Range("E9").Value = 0.008403361344537815
$ws.Range("F9").Value = 0.07563025210084033
$ws.Range("J9").Value = 0.1428571428571428
$ws.Range("Q9").Value = 0.1596638655462185
$ws.Range("R9").Value = 0.06722689075630252
$ws.Range("S9").Value = 0.4033613445378151
$ws.Range("B10").Value = 0.09013605442176871
$ws.Range("D10").Value = 0.02295918367346939
$ws.Range("E10").Value = 0.0008503401360544217
$ws.Range("F10").Value = 0.07142857142857142
$ws.Range("J10").Value = 0.1335034013605442
$ws.Range("O10").Value = 0.01530612244897959
$ws.Range("Q10").Value = 0.2304421768707483
$ws.Range("R10").Value = 0.07653061224489796
$ws.Range("S10").Value = 0.358843537414966
$ws.Range("G11").Value = 0.1524163568773234
$ws.Range("J11").Value = 0.09665427509293681
$ws.Range("K11").Value = 0.2230483271375465
$ws.Range("L11").Value = 0.5241635687732342
$ws.Range("S11").Value = 0.003717472118959108
$ws.Range("G12").Value = 0.7622377622377622
$ws.Range("J12").Value = 0.2027972027972028
$ws.Range("K12").Value = 0.006993006993006993
$ws.Range("L12").Value = 0.01398601398601399
$ws.Range("S12").Value = 0.01398601398601399
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.3333333333333333
$ws.Range("F15").Value = 0.02702702702702703
$ws.Range("H15").Value = 0.1405405405405405
$ws.Range("I15").Value = 0.04324324324324325
$ws.Range("J15").Value = 0.4
$ws.Range("K15").Value = 0.1135135135135135
$ws.Range("M15").Value = 0.01621621621621622
$ws.Range("O15").Value = 0.06486486486486487
$ws.Range("S15").Value = 0.1945945945945946
$ws.Range("F16").Value = 0.03684210526315789
$ws.Range("H16").Value = 0.1526315789473684
$ws.Range("I16").Value = 0.07368421052631578
$ws.Range("J16").Value = 0.4421052631578947
$ws.Range("K16").Value = 0.09473684210526316
$ws.Range("M16").Value = 0.01578947368421053
$ws.Range("O16").Value = 0.06315789473684211
$ws.Range("S16").Value = 0.1210526315789474
$ws.Range("F17").Value = 0.03044496487119438
$ws.Range("H17").Value = 0.1569086651053864
$ws.Range("I17").Value = 0.07259953161592506
$ws.Range("J17").Value = 0.4426229508196721
$ws.Range("K17").Value = 0.1124121779859485
$ws.Range("M17").Value = 0.01873536299765808
$ws.Range("O17").Value = 0.04449648711943794
$ws.Range("S17").Value = 0.1217798594847775
$ws.Range("F18").Value = 0.006060606060606061
$ws.Range("H18").Value = 0.2
$ws.Range("I18").Value = 0.0303030303030303
$ws.Range("J18").Value = 0.4363636363636363
$ws.Range("K18").Value = 0.05454545454545454
$ws.Range("M18").Value = 0.0303030303030303
$ws.Range("O18").Value = 0.07878787878787878
$ws.Range("S18").Value = 0.1636363636363636
$ws.Range("F19").Value = 0.01923076923076923
$ws.Range("H19").Value = 0.2344322344322344
$ws.Range("I19").Value = 0.05494505494505494
$ws.Range("J19").Value = 0.3663003663003663
$ws.Range("K19").Value = 0.108974358974359
$ws.Range("M19").Value = 0.02564102564102564
$ws.Range("O19").Value = 0.07326007326007326
$ws.Range("S19").Value = 0.1172161172161172

Write-Output "Applied 107 cell updates"
